$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New other_party_org rows for "Сміляньска/Смілянська міська організація
# Партії Зелених України" across the 4 quarters of 2021 (two EDRPOU-code
# variants per quarter).
$rows = @(
    @("Сміляньска міська організація Партії Зелених України", "36463324", "1 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "9032f690-0182-11ef-938a-5514903835ea"),
    @("Сміляньска міська організація Партії Зелених України", "36463324", "1 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "9032f690-0182-11ef-938a-5514903835ea"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "1 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "9032f690-0182-11ef-938a-5514903835ea"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "1 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "9032f690-0182-11ef-938a-5514903835ea"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "2 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "a38ae5e0-0182-11ef-95c4-e7bfad33aba2"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "2 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "a38ae5e0-0182-11ef-95c4-e7bfad33aba2"),
    @("Сміляньска міська організація Партії Зелених України", "36463324", "2 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "a38ae5e0-0182-11ef-95c4-e7bfad33aba2"),
    @("Сміляньска міська організація Партії Зелених України", "36463324", "2 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "a38ae5e0-0182-11ef-95c4-e7bfad33aba2"),
    @("Сміляньска міська організація Партії Зелених України", "36463324", "3 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "b2cb7830-0182-11ef-ae0a-87878d7061af"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "3 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "b2cb7830-0182-11ef-ae0a-87878d7061af"),
    @("Сміляньска міська організація Партії Зелених України", "36463324", "4 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "c5f838d0-0182-11ef-95c4-e7bfad33aba2"),
    @("Смілянська міська організація Партії Зелених України", "26043372", "4 квартал", "2021", "ПАРТІЯ ЗЕЛЕНИХ УКРАЇНИ", "00047728", "c5f838d0-0182-11ef-95c4-e7bfad33aba2")
)

$startRow = 76
$endRow = $startRow + $rows.Length - 1

# Columns B (EDRPOU code), D (year) and F (party EDRPOU code) hold values
# that look numeric ("36463324", "2021", "00047728"); mark each column
# range as Text up front so the values are stored verbatim (keeping the
# leading zero in "00047728") instead of being auto-coerced to numbers,
# then strip the temporary formatting back off so the cells end up
# plain/unstyled, matching the rest of the sheet.
$textColRanges = @(
    $ws.Range("B$startRow`:B$endRow"),
    $ws.Range("D$startRow`:D$endRow"),
    $ws.Range("F$startRow`:F$endRow")
)
foreach ($colRange in $textColRanges) {
    $colRange.NumberFormat = "@"
}

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

foreach ($colRange in $textColRanges) {
    $colRange.ClearFormats()
}
